$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Intro paragraph: replace "Section 2 Team 1" with the team name "Xeno".
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Section 2 Team 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $introTeam = $d.Range($find.Start, $find.End)
    $introTeam.Text = "Xeno"
}

# ---------------------------------------------------------------------------
# 2. "Team: ____" paragraph -> bold "Team: Xeno" line.
# ---------------------------------------------------------------------------
$teamPara = $d.Paragraphs(3)
$teamRange = $d.Range($teamPara.Range.Start, $teamPara.Range.End)
$teamRange.Text = "Team: Xeno"
$teamRange2 = $d.Range($teamPara.Range.Start, $teamPara.Range.End)
$teamRange2.Bold = 1

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark paragraph so it sits right after the Team
#    paragraph (i.e. right before the "Name: Edward LaFemina" block) instead
#    of right before the final blank paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$afterTeam = $d.Paragraphs(4)
$newBookmarkRange = $afterTeam.Range
$newBookmarkRange.Collapse(1)
$newBookmarkRange.Bookmarks.Add("_GoBack")

# ---------------------------------------------------------------------------
# 4. Simplify Edward LaFemina's "Date:" line from three runs into one.
# ---------------------------------------------------------------------------
$datePara = $d.Paragraphs(7)
$dateRange = $d.Range($datePara.Range.Start, $datePara.Range.End)
$dateRange.Text = "Date: __________________________"
